$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of infaunal sample data (rows 140-180), columns A (block), B (number), C (taxon), D (abund), E (notes)
$data = New-Object 'object[][]' 41

$data[0]  = @('A', 6,  'lottia_paradigitalis',   3,  $null)
$data[1]  = @('A', 6,  'mytilus',                1,  $null)
$data[2]  = @('A', 6,  'amphipoda_1',             2,  $null)
$data[3]  = @('A', 6,  'littorina_scutulata',     1,  $null)
$data[4]  = @('C', 16, 'cirripedia_cyprid',       3,  $null)
$data[5]  = @('C', 16, 'mytilus',                 2,  $null)
$data[6]  = @('C', 16, 'amphipoda_1',             3,  $null)
$data[7]  = @('C', 16, 'lottia_paradigitalis',    2,  $null)
$data[8]  = @('C', 16, 'littorina_scutulata',     11, $null)
$data[9]  = @('C', 16, 'littorina_sitkana',       81, $null)
$data[10] = @('D', 15, 'cirripedia_cyprid',       4,  $null)
$data[11] = @('D', 15, 'mytilus',                 6,  $null)
$data[12] = @('D', 15, 'copepoda',                1,  $null)
$data[13] = @('D', 15, 'emplectonema_gracile',    1,  $null)
$data[14] = @('D', 15, 'lottia_paradigitalis',    5,  $null)
$data[15] = @('D', 15, 'pagurus_hirsutiusculus',  1,  $null)
$data[16] = @('D', 15, 'littorina_scutulata',     59, $null)
$data[17] = @('D', 15, 'isopoda_1',               3,  $null)
$data[18] = @('D', 15, 'littorina_sitkana',       11, $null)
$data[19] = @('D', 15, 'worm_thing',              1,  $null)
$data[20] = @('D', 15, 'amphipoda_1',             29, $null)
$data[21] = @('D', 15, 'limpet_recruit',          1,  $null)
$data[22] = @('F', 4,  'cirripedia_cyprid',       7,  $null)
$data[23] = @('F', 4,  'mytilus',                 15, $null)
$data[24] = @('F', 4,  'lottia_paradigitalis',    3,  $null)
$data[25] = @('F', 4,  'emplectonema_gracile',    1,  $null)
$data[26] = @('F', 4,  'worm_thing',              2,  $null)
$data[27] = @('F', 4,  'isopoda_1',               3,  $null)
$data[28] = @('F', 4,  'amphipoda_1',             26, $null)
$data[29] = @('F', 4,  'littorina_scutulata',     57, $null)
$data[30] = @('F', 4,  'littorina_sitkana',       76, $null)
$data[31] = @('F', 4,  'bivalvia_1',              1,  $null)
$data[32] = @('D', 6,  'cirripedia_cyprid',       10, $null)
$data[33] = @('D', 6,  'mytilus',                 7,  $null)
$data[34] = @('D', 6,  'polychaeta_2',            2,  'grey, short')
$data[35] = @('D', 6,  'lottia_paradigitalis',    2,  $null)
$data[36] = @('D', 6,  'amphipoda_1',             14, $null)
$data[37] = @('D', 6,  'littorina_scutulata',     55, $null)
$data[38] = @('D', 6,  'littorina_sitkana',       51, $null)
$data[39] = @('D', 6,  'oedoperna_larvae',        2,  $null)
$data[40] = @('D', 6,  'insecta_2',               1,  'tick')

# Rows where the note (column E) was authored before the taxon (column C),
# so the shared-string table records the note's text ahead of the taxon's.
$notesBeforeTaxon = @(34)

$startRow = 140
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($notesBeforeTaxon -contains $i) {
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
    } else {
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        if ($row[4] -ne $null) {
            $ws.Cells.Item($r, 5).Value = $row[4]
        }
    }
}

# Update active cell / selection to match the saved workbook state
$ws.Range("G176").Select()
